# Fix error in oracle: the "ERROR" check formulas in column S were only
# testing the raw 1..5 rating value (e.g. B8>3) instead of the weighted
# score (rating * weight, e.g. B8*C8>3). Update the three shared-formula
# "master" cells (S8, S9, S73) so the fix propagates through the shared
# formula ranges S9:S72 and S73:S105, which in turn recalculates the
# downstream U and V columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BSN-Test Set")

$ws.Range("S8").Formula = '=IF(IF(B8*C8>3,1,0)+IF(D8*E8>3,1,0)+IF(F8*G8>3,1,0)+IF(H8*I8>3,1,0)+IF(J8*K8>3,1,0)+IF(L8*M8>3,1,0)>0,"ERROR",0)'
$ws.Range("S9:S72").Formula = '=IF(IF(B9*C9>3,1,0)+IF(D9*E9>3,1,0)+IF(F9*G9>3,1,0)+IF(H9*I9>3,1,0)+IF(J9*K9>3,1,0)+IF(L9*M9>3,1,0)>0,"ERROR",0)'
$ws.Range("S73:S105").Formula = '=IF(IF(B73*C73>3,1,0)+IF(D73*E73>3,1,0)+IF(F73*G73>3,1,0)+IF(H73*I73>3,1,0)+IF(J73*K73>3,1,0)+IF(L73*M73>3,1,0)>0,"ERROR",0)'

# Update the saved view state to match (scrolled/selected cell changed
# incidentally as part of the author's editing session).
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("S8").Select()
